$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last row that actually has data in column A (Beteckning),
# using the classic "Ctrl+Up from the bottom" idiom (xlUp = -4162).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

if ($lastRow -lt 2) {
    # Fallback in case the above did not resolve as expected.
    $lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1
}
if ($lastRow -lt 2) {
    $lastRow = 259
}

# Column C ("Förändrad") holds the last-changed date for every logging
# notice row. Bump it from 2026-02-21 (serial 46074) to 2026-02-22
# (serial 46075) for every data row, leaving everything else untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq 46074) {
        $cell.Value2 = 46075
    }
}
